$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.703.66'
$ws.Range("E2").Value = '  +5.56%  '

$ws.Range("D3").Value = '3.329.35'
$ws.Range("E3").Value = '  +4.98%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '553.89'
$ws.Range("E5").Value = '  +3.81%  '

$ws.Range("D6").Value = '151.61'
$ws.Range("E6").Value = '  +5.27%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  +1.40%  '

$ws.Range("D9").Value = '7.49'
$ws.Range("E9").Value = '  +3.04%  '

$ws.Range("E10").Value = '  +4.73%  '

$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").Value = '3.900.78'
$ws.Range("E12").Value = '  +4.91%  '

$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("E14").Value = '  +5.47%  '

$ws.Range("D15").Value = '26.97'
$ws.Range("E15").Value = '  +4.10%  '

$ws.Range("D16").Value = '62.661.26'
$ws.Range("E16").Value = '  +5.43%  '

$ws.Range("D17").Value = '3.332.10'
$ws.Range("E17").Value = '  +4.85%  '

$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  +5.19%  '

$ws.Range("D19").Value = '13.82'
$ws.Range("E19").Value = '  +6.72%  '

$ws.Range("D20").Value = '8.49'
$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").Value = '385.65'
$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("D22").Value = '1.01'
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("E23").Value = '  +2.04%  '

$ws.Range("D24").Value = '70.94'
$ws.Range("E24").Value = '  +1.16%  '

$ws.Range("E25").Value = '  +4.50%  '

$ws.Range("D26").Value = '8.83'
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = '0.0₃0972'
$ws.Range("E27").Value = '  +8.49%  '

$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("D29").Value = '6.46'
$ws.Range("E29").Value = '  +5.35%  '

$ws.Range("D30").Value = '1.98'
$ws.Range("E30").Value = '  +4.07%  '

$ws.Range("D31").Value = '22.97'
$ws.Range("E31").Value = '  +3.11%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '5.58'
$ws.Range("E32").Value = '  +5.63%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  +11.37%  '

$ws.Range("D34").Value = '6.74'
$ws.Range("E34").Value = '  +4.49%  '

$ws.Range("E35").Value = '  +11.30%  '

$ws.Range("D36").Value = '160.48'
$ws.Range("E36").Value = '  +2.56%  '

$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  +13.25%  '

$ws.Range("D38").Value = '27.05'
$ws.Range("E38").Value = '  +7.10%  '

$ws.Range("D39").Value = '2.856.96'
$ws.Range("E39").Value = '  +4.31%  '

$ws.Range("E40").Value = '  +3.77%  '

$ws.Range("D41").Value = '0.0315'
$ws.Range("E41").Value = '  +8.81%  '

$ws.Range("D42").Value = '4.34'
$ws.Range("E42").Value = '  +1.46%  '

$ws.Range("D43").Value = '0.752'
$ws.Range("E43").Value = '  +4.07%  '

$ws.Range("D44").Value = '40.66'
$ws.Range("E44").Value = '  +3.37%  '

$ws.Range("E45").Value = '  +4.07%  '

$ws.Range("D46").Value = '22.07'
$ws.Range("E46").Value = '  +8.14%  '

$ws.Range("D47").Value = '3.372.96'
$ws.Range("E47").Value = '  +4.91%  '

$ws.Range("E48").Value = '  +4.32%  '

$ws.Range("D49").Value = '6.30'
$ws.Range("E49").Value = '  +2.20%  '

$ws.Range("D50").Value = '0.811'
$ws.Range("E50").Value = '  +5.69%  '

$ws.Range("D51").Value = '283.16'
$ws.Range("E51").Value = '  +8.70%  '

